$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 described spare parts (B12 "Project box", C12 "none", G12 "Onwon ",
# H12 "Amazon", I12 hyperlink to "https://amzn.to/2CqpfJ8 ") is removed from
# the BOM; only I12 keeps its (Hyperlink) cell style with no value.
$ws.Range("A12:H12").ClearContents()
$ws.Range("I12").ClearContents()

# Drop the now-orphaned hyperlink on I12, but keep the other six intact.
# (This engine's Hyperlinks.Delete() clears the whole sheet collection, so
# rebuild the ones that should survive.)
$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("I2"), "https://amzn.to/32xisrJ")
$null = $ws.Hyperlinks.Add($ws.Range("I3"), "https://amzn.to/3jiLVvC")
$null = $ws.Hyperlinks.Add($ws.Range("I8"), "https://amzn.to/30r0jJy")
$null = $ws.Hyperlinks.Add($ws.Range("I9"), "https://amzn.to/3fCVfrP")
$null = $ws.Hyperlinks.Add($ws.Range("I10"), "https://amzn.to/30qN8Il")
$null = $ws.Hyperlinks.Add($ws.Range("I11"), "https://amzn.to/32wXrgN")

# Re-adding hyperlinks resets the cell style to a fresh copy of the
# "Hyperlink" style; restore the original shared style on each cell so the
# stylesheet doesn't pick up a duplicate.
$ws.Range("I2").Style = "Hyperlink"
$ws.Range("I3").Style = "Hyperlink"
$ws.Range("I8").Style = "Hyperlink"
$ws.Range("I9").Style = "Hyperlink"
$ws.Range("I10").Style = "Hyperlink"
$ws.Range("I11").Style = "Hyperlink"

# Update the saved selection/view state to the new used range.
$ws.Range("A1:I11").Select()
